$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (duplicate-renumbered issues -> bumped revision suffix)
$ws.Name = "EvaChecks (6)"

# Row 112 <- was the "Overlapping CE Participation Records" row (old row 116)
$ws.Range("B112").Value = "High Priority"
$ws.Range("C112").Value = "Overlapping CE Participation Records"
$ws.Range("D112").Value = "2.09 CE Participation Status"
$ws.Range("E112").Value = "This project has more than one CE Participation record that covers the same time period. Please be sure you are ending any records that are no longer accurate before creating a new record."
$ws.Range("G112").Value = 128

# Row 113 <- was the "Overlapping HMIS Participation Records" row (old row 117), renumbered 129 -> 131
$ws.Range("B113").Value = "High Priority"
$ws.Range("C113").Value = "Overlapping HMIS Participation Records"
$ws.Range("D113").Value = "2.08 HMIS Participation"
$ws.Range("E113").Value = "This project has more than one HMIS Participation record that covers the same time period. Please be sure you are ending any records that are no longer accurate before creating a new record."
$ws.Range("G113").Value = 131

# Row 114 <- was "Missing Address" (old row 112)
$ws.Range("B114").Value = "Error"
$ws.Range("C114").Value = "Missing Address"
$ws.Range("D114").Value = "2.03.1-2.03.7 - Continuum of Care Information"
$ws.Range("E114").Value = "Please ensure geography information for projects is complete."
$ws.Range("G114").Value = 42

# Row 115 <- was "Missing RRH SubType" (old row 113)
$ws.Range("B115").Value = "Error"
$ws.Range("C115").Value = "Missing RRH SubType"
$ws.Range("D115").Value = "2.02A RRH SubType, 2.02.6 Project Type"
$ws.Range("E115").Value = "All RRH projects must have an RRH SubType. Please update the data at the project level."
$ws.Range("G115").Value = 110

# Row 116 <- was "No Inventory Records" (old row 114)
$ws.Range("B116").Value = "Error"
$ws.Range("C116").Value = "No Inventory Records"
$ws.Range("D116").Value = "2.07 - Bed and Unit Inventory Information, 2.02.6 - Project Type"
$ws.Range("E116").Value = "Residential projects should have inventory data. Please enter inventory in HMIS for the project."
$ws.Range("G116").Value = 43

# Row 117 <- was "Operating End precedes Inventory End" (old row 115)
$ws.Range("B117").Value = "Error"
$ws.Range("C117").Value = "Operating End precedes Inventory End"
$ws.Range("D117").Value = "2.02.3 - Project Information: Operating Start Date, 3.10 - Project Information: Project Entry Date"
$ws.Range("E117").Value = "Inventory records should end whenever a project stops operating. Please correct either the inventory dates or the Project Operating End Date."
$ws.Range("G117").Value = 44

# Row 118 stays "RRH-SO has active bed inventory" but renumbered 130 -> 132
$ws.Range("G118").Value = 132

# Scroll/selection state saved with the workbook
$ws.Application.ActiveWindow.ScrollRow = 95
$ws.Range("C118").Select()
